$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.868.02"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "3.168.59"
$ws.Range("E3").Value = "  +4.14%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.74"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.11"
$ws.Range("E6").Value = "  +6.74%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.168.01"
$ws.Range("E8").Value = "  +4.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("E10").Value = "  +6.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.25"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.501"
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("E13").Value = "  +17.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.68"
$ws.Range("E14").Value = "  +6.02%  "
$ws.Range("D15").Value = "3.688.27"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "64.985.48"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "3.170.30"
$ws.Range("E17").Value = "  +4.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.16"
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.61"
$ws.Range("E20").Value = "  +7.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.82"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("E22").Value = "  +6.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.32"
$ws.Range("E23").Value = "  +5.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.79"
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.25"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +5.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.99"
$ws.Range("E28").Value = "  +11.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  +7.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.81"
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.75"
$ws.Range("E31").Value = "  +13.20%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +4.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.32"
$ws.Range("E34").Value = "  +10.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("E35").Value = "  +6.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.72"
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0891"
$ws.Range("E37").Value = "  +10.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "472.42"
$ws.Range("E38").Value = "  +7.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  +12.06%  "
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.64"
$ws.Range("E41").Value = "  +4.78%  "
$ws.Range("D42").Value = "3.064.41"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.284"
$ws.Range("E44").Value = "  +6.14%  "
$ws.Range("E45").Value = "  +8.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.14"
$ws.Range("E46").Value = "  +4.96%  "
$ws.Range("D47").Value = "0.0₃0606"
$ws.Range("E47").Value = "  +18.73%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  +8.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.22"
$ws.Range("E51").Value = "  +1.77%  "
